$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 15; existing rows 15.. shift down to 16..
$ws.Rows("15:15").Insert()

# Populate the new row 15 with data (same pattern/columns as surrounding
# rows for this market/product, with its own date/volume/price figures).
$ws.Range("A15").Value2 = 2
$ws.Range("B15").Value2 = "Comercializadora del Agro de Limarí"
$ws.Range("C15").Value2 = "Coquimbo"
$ws.Range("D15").Value2 = 44560
$ws.Range("E15").Value2 = 4
$ws.Range("F15").Value2 = 100112030
$ws.Range("G15").Value2 = "Poroto granado"
$ws.Range("H15").Value2 = "Sin especificar"
$ws.Range("I15").Value2 = "Primera"
$ws.Range("J15").Value2 = 680
$ws.Range("K15").Value2 = 16000
$ws.Range("L15").Value2 = 17000
$ws.Range("M15").Value2 = 16500
$ws.Range("N15").Value2 = "$/malla 25 kilos"
$ws.Range("O15").Value2 = "Provincia de Limarí"
$ws.Range("P15").Value2 = 660
$ws.Range("Q15").Value2 = 25
$ws.Range("R15").Value2 = "Hortaliza"
